$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp shared string (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 17:35"

# Update province statistics
$ws.Range("B4").Value = 66860
$ws.Range("D4").Value = 57948
$ws.Range("E4").Value = 8912

$ws.Range("B5").Value = 55888
$ws.Range("D5").Value = 49867
$ws.Range("E5").Value = 6021

$ws.Range("B6").Value = 18586
$ws.Range("D6").Value = 16626

$ws.Range("B7").Value = 16739
$ws.Range("D7").Value = 13826
$ws.Range("E7").Value = 2913

$ws.Range("B9").Value = 12502
$ws.Range("D9").Value = 11131
$ws.Range("E9").Value = 1371

$ws.Range("B14").Value = 5551
$ws.Range("D14").Value = 4708

$ws.Range("B16").Value = 5136
$ws.Range("D16").Value = 4633

$ws.Range("B20").Value = 4033
$ws.Range("D20").Value = 3680
$ws.Range("E20").Value = 353

$ws.Range("B32").Value = 2374
$ws.Range("D32").Value = 2070
$ws.Range("E32").Value = 304

$ws.Range("B33").Value = 2300
$ws.Range("D33").Value = 2145
